$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# short-url column (B) - shared string changes for every row, set via B2 (shared string)
$ws.Range("B2:B40").Value = "6qjL6E"

# The source cells store purely-numeric-looking values as TEXT (shared
# strings) rather than numbers. A direct Range.Value = "123..." assignment
# gets auto-coerced to a numeric cell (losing the text type and bumping the
# cell style to a new quote-prefixed xf). To keep the text type AND the
# original style, round-trip the new value through a formula that evaluates
# to text, then paste-special just the value on top of the target cell.
$helper = $ws.Range("ZZ1")

function Set-TextValue($rangeAddress, $textValue) {
    $helper.Formula = '="' + $textValue + '"'
    $helper.Copy()
    $ws.Range($rangeAddress).PasteSpecial(-4163)  # xlPasteValues
}

# 2024 row (row 39): idps (Q) and returned_idps (R)
Set-TextValue "Q39" "3513867"
Set-TextValue "R39" "377566"

# 2024 row (row 40): stateless (S)
Set-TextValue "S40" "619429"

$helper.ClearContents()
$excel.CutCopyMode = $false
